# The workbook tracks daily "Brócoli" price observations in a single
# growing table (Sheet1). A new observation was inserted as row 352,
# pushing every subsequent row (old 352..442) down by one (new 353..443).
# This mirrors what happens in Excel when a user selects row 352 and
# inserts a new blank row above it, then fills in the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 352; existing rows 352-442 shift to 353-443.
$ws.Rows("352:352").Insert()

# Populate the newly inserted row 352 with the new observation.
$ws.Range("A352").Value = 10
$ws.Range("B352").Value = "Vega Modelo de Temuco"
$ws.Range("C352").Value = "La Araucanía"
$ws.Range("D352").Value = 44736
$ws.Range("E352").Value = 9
$ws.Range("F352").Value = 100112023
$ws.Range("G352").Value = "Brócoli"
$ws.Range("H352").Value = "Sin especificar"
$ws.Range("I352").Value = "Primera"
$ws.Range("J352").Value = 400
$ws.Range("K352").Value = 1000
$ws.Range("L352").Value = 1200
$ws.Range("M352").Value = 1100
$ws.Range("N352").Value = '$/unidad'
$ws.Range("O352").Value = "Región Metropolitana"
$ws.Range("P352").Value = 1100
$ws.Range("Q352").Value = 1
$ws.Range("R352").Value = "Hortaliza"
